$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 27: remove the now-unused D27/E27 cells, label F27, keep the G27 parameters.json formula
$ws.Range("D27").Clear()
$ws.Range("E27").Clear()

$ws.Range("F26").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("F27").Value = "parameters.json:"

# Row 29: add the "presets HTML:" label and the two <option> builder formulas
$ws.Range("F30").Copy()
$ws.Range("F29:I29").PasteSpecial(-4122)

$ws.Range("F29").Value = "presets HTML:"
$ws.Range("G29").Formula = '="<option value=""" & D$2 & """>" & D$2 & "</option>"'
$ws.Range("H29").Formula = '="<option value=""" & E$2 & """>" & E$2 & "</option>"'

# Match the author's final selection
$ws.Range("G30").Select()
